$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2023-12-19 Tuesday" "2023-12-20 Wednesday"

Replace-Text "92×24=" "47×82="
Replace-Text "92×98=" "12×70="
Replace-Text "86×23=" "70×51="
Replace-Text "34×55=" "33×26="
Replace-Text "85×92=" "43×94="

Replace-Text "99×96=" "19×30="
Replace-Text "22×57=" "56×67="
Replace-Text "80×46=" "83×61="
Replace-Text "34×61=" "12×20="
Replace-Text "28×96=" "20×70="

Replace-Text "51×37=" "67×42="
Replace-Text "87×25=" "45×77="
Replace-Text "72×59=" "25×80="
Replace-Text "19×65=" "76×34="
Replace-Text "46×36=" "75×51="

Replace-Text "64×99=" "26×93="
Replace-Text "12×87=" "45×57="
Replace-Text "65×76=" "58×13="
Replace-Text "34×64=" "84×18="
Replace-Text "54×60=" "63×28="

Replace-Text "29×16=" "50×56="
Replace-Text "42×79=" "43×34="
Replace-Text "77×20=" "71×48="
Replace-Text "80×35=" "35×76="
Replace-Text "69×89=" "72×83="
